# Shift production-prediction data forward by 8 days and update the
# Prediction (EET-adjusted) values, per the "Adjusting to EET" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Treat column A as plain text first, so the date-like strings are not
# auto-converted into date serial numbers by Excel.
$ws.Range("A2:A170").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2024-04-03"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 1).Value = "2024-04-03"
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(4, 1).Value = "2024-04-03"
$ws.Cells.Item(4, 3).Value = 0.066
$ws.Cells.Item(5, 1).Value = "2024-04-03"
$ws.Cells.Item(5, 3).Value = 0.013
$ws.Cells.Item(6, 1).Value = "2024-04-03"
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = "2024-04-03"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = "2024-04-03"
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = "2024-04-03"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = "2024-04-04"
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = "2024-04-04"
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = "2024-04-04"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 1).Value = "2024-04-04"
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 1).Value = "2024-04-04"
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(15, 1).Value = "2024-04-04"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(16, 1).Value = "2024-04-04"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(17, 1).Value = "2024-04-04"
$ws.Cells.Item(17, 3).Value = 0.012
$ws.Cells.Item(18, 1).Value = "2024-04-04"
$ws.Cells.Item(18, 3).Value = 0.046
$ws.Cells.Item(19, 1).Value = "2024-04-04"
$ws.Cells.Item(19, 3).Value = 0.117
$ws.Cells.Item(20, 1).Value = "2024-04-04"
$ws.Cells.Item(20, 3).Value = 0.218
$ws.Cells.Item(21, 1).Value = "2024-04-04"
$ws.Cells.Item(21, 3).Value = 0.271
$ws.Cells.Item(22, 1).Value = "2024-04-04"
$ws.Cells.Item(22, 3).Value = 0.305
$ws.Cells.Item(23, 1).Value = "2024-04-04"
$ws.Cells.Item(23, 3).Value = 0.326
$ws.Cells.Item(24, 1).Value = "2024-04-04"
$ws.Cells.Item(24, 3).Value = 0.319
$ws.Cells.Item(25, 1).Value = "2024-04-04"
$ws.Cells.Item(25, 3).Value = 0.307
$ws.Cells.Item(26, 1).Value = "2024-04-04"
$ws.Cells.Item(26, 3).Value = 0.241
$ws.Cells.Item(27, 1).Value = "2024-04-04"
$ws.Cells.Item(27, 3).Value = 0.185
$ws.Cells.Item(28, 1).Value = "2024-04-04"
$ws.Cells.Item(28, 3).Value = 0.082
$ws.Cells.Item(29, 1).Value = "2024-04-04"
$ws.Cells.Item(29, 3).Value = 0.016
$ws.Cells.Item(30, 1).Value = "2024-04-04"
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 1).Value = "2024-04-04"
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 1).Value = "2024-04-04"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 1).Value = "2024-04-04"
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 1).Value = "2024-04-05"
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 1).Value = "2024-04-05"
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 1).Value = "2024-04-05"
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 1).Value = "2024-04-05"
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(38, 1).Value = "2024-04-05"
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 1).Value = "2024-04-05"
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 1).Value = "2024-04-05"
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 1).Value = "2024-04-05"
$ws.Cells.Item(41, 3).Value = 0.014
$ws.Cells.Item(42, 1).Value = "2024-04-05"
$ws.Cells.Item(42, 3).Value = 0.061
$ws.Cells.Item(43, 1).Value = "2024-04-05"
$ws.Cells.Item(43, 3).Value = 0.17
$ws.Cells.Item(44, 1).Value = "2024-04-05"
$ws.Cells.Item(44, 3).Value = 0.252
$ws.Cells.Item(45, 1).Value = "2024-04-05"
$ws.Cells.Item(45, 3).Value = 0.345
$ws.Cells.Item(46, 1).Value = "2024-04-05"
$ws.Cells.Item(46, 3).Value = 0.404
$ws.Cells.Item(47, 1).Value = "2024-04-05"
$ws.Cells.Item(47, 3).Value = 0.411
$ws.Cells.Item(48, 1).Value = "2024-04-05"
$ws.Cells.Item(48, 3).Value = 0.417
$ws.Cells.Item(49, 1).Value = "2024-04-05"
$ws.Cells.Item(49, 3).Value = 0.356
$ws.Cells.Item(50, 1).Value = "2024-04-05"
$ws.Cells.Item(50, 3).Value = 0.287
$ws.Cells.Item(51, 1).Value = "2024-04-05"
$ws.Cells.Item(51, 3).Value = 0.189
$ws.Cells.Item(52, 1).Value = "2024-04-05"
$ws.Cells.Item(52, 3).Value = 0.083
$ws.Cells.Item(53, 1).Value = "2024-04-05"
$ws.Cells.Item(53, 3).Value = 0.015
$ws.Cells.Item(54, 1).Value = "2024-04-05"
$ws.Cells.Item(54, 3).Value = 0
$ws.Cells.Item(55, 1).Value = "2024-04-05"
$ws.Cells.Item(55, 3).Value = 0
$ws.Cells.Item(56, 1).Value = "2024-04-05"
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(57, 1).Value = "2024-04-05"
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(58, 1).Value = "2024-04-06"
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(59, 1).Value = "2024-04-06"
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(60, 1).Value = "2024-04-06"
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(61, 1).Value = "2024-04-06"
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(62, 1).Value = "2024-04-06"
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(63, 1).Value = "2024-04-06"
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(64, 1).Value = "2024-04-06"
$ws.Cells.Item(64, 3).Value = 0
$ws.Cells.Item(65, 1).Value = "2024-04-06"
$ws.Cells.Item(65, 3).Value = 0.015
$ws.Cells.Item(66, 1).Value = "2024-04-06"
$ws.Cells.Item(66, 3).Value = 0.095
$ws.Cells.Item(67, 1).Value = "2024-04-06"
$ws.Cells.Item(67, 3).Value = 0.217
$ws.Cells.Item(68, 1).Value = "2024-04-06"
$ws.Cells.Item(68, 3).Value = 0.352
$ws.Cells.Item(69, 1).Value = "2024-04-06"
$ws.Cells.Item(69, 3).Value = 0.476
$ws.Cells.Item(70, 1).Value = "2024-04-06"
$ws.Cells.Item(70, 3).Value = 0.525
$ws.Cells.Item(71, 1).Value = "2024-04-06"
$ws.Cells.Item(71, 3).Value = 0.5
$ws.Cells.Item(72, 1).Value = "2024-04-06"
$ws.Cells.Item(72, 3).Value = 0.475
$ws.Cells.Item(73, 1).Value = "2024-04-06"
$ws.Cells.Item(73, 3).Value = 0.404
$ws.Cells.Item(74, 1).Value = "2024-04-06"
$ws.Cells.Item(74, 3).Value = 0.283
$ws.Cells.Item(75, 1).Value = "2024-04-06"
$ws.Cells.Item(75, 3).Value = 0.193
$ws.Cells.Item(76, 1).Value = "2024-04-06"
$ws.Cells.Item(76, 3).Value = 0.085
$ws.Cells.Item(77, 1).Value = "2024-04-06"
$ws.Cells.Item(77, 3).Value = 0.015
$ws.Cells.Item(78, 1).Value = "2024-04-06"
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(79, 1).Value = "2024-04-06"
$ws.Cells.Item(79, 3).Value = 0
$ws.Cells.Item(80, 1).Value = "2024-04-06"
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(81, 1).Value = "2024-04-06"
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(82, 1).Value = "2024-04-07"
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(83, 1).Value = "2024-04-07"
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(84, 1).Value = "2024-04-07"
$ws.Cells.Item(84, 3).Value = 0
$ws.Cells.Item(85, 1).Value = "2024-04-07"
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(86, 1).Value = "2024-04-07"
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(87, 1).Value = "2024-04-07"
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(88, 1).Value = "2024-04-07"
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(89, 1).Value = "2024-04-07"
$ws.Cells.Item(89, 3).Value = 0.024
$ws.Cells.Item(90, 1).Value = "2024-04-07"
$ws.Cells.Item(90, 3).Value = 0.118
$ws.Cells.Item(91, 1).Value = "2024-04-07"
$ws.Cells.Item(91, 3).Value = 0.256
$ws.Cells.Item(92, 1).Value = "2024-04-07"
$ws.Cells.Item(92, 3).Value = 0.385
$ws.Cells.Item(93, 1).Value = "2024-04-07"
$ws.Cells.Item(93, 3).Value = 0.483
$ws.Cells.Item(94, 1).Value = "2024-04-07"
$ws.Cells.Item(94, 3).Value = 0.544
$ws.Cells.Item(95, 1).Value = "2024-04-07"
$ws.Cells.Item(95, 3).Value = 0.558
$ws.Cells.Item(96, 1).Value = "2024-04-07"
$ws.Cells.Item(96, 3).Value = 0.531
$ws.Cells.Item(97, 1).Value = "2024-04-07"
$ws.Cells.Item(97, 3).Value = 0.455
$ws.Cells.Item(98, 1).Value = "2024-04-07"
$ws.Cells.Item(98, 3).Value = 0.346
$ws.Cells.Item(99, 1).Value = "2024-04-07"
$ws.Cells.Item(99, 3).Value = 0.242
$ws.Cells.Item(100, 1).Value = "2024-04-07"
$ws.Cells.Item(100, 3).Value = 0.124
$ws.Cells.Item(101, 1).Value = "2024-04-07"
$ws.Cells.Item(101, 3).Value = 0.024
$ws.Cells.Item(102, 1).Value = "2024-04-07"
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(103, 1).Value = "2024-04-07"
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(104, 1).Value = "2024-04-07"
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(105, 1).Value = "2024-04-07"
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(106, 1).Value = "2024-04-08"
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(107, 1).Value = "2024-04-08"
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(108, 1).Value = "2024-04-08"
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(109, 1).Value = "2024-04-08"
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(110, 1).Value = "2024-04-08"
$ws.Cells.Item(110, 3).Value = 0
$ws.Cells.Item(111, 1).Value = "2024-04-08"
$ws.Cells.Item(111, 3).Value = 0
$ws.Cells.Item(112, 1).Value = "2024-04-08"
$ws.Cells.Item(112, 3).Value = 0
$ws.Cells.Item(113, 1).Value = "2024-04-08"
$ws.Cells.Item(113, 3).Value = 0.022
$ws.Cells.Item(114, 1).Value = "2024-04-08"
$ws.Cells.Item(114, 3).Value = 0.118
$ws.Cells.Item(115, 1).Value = "2024-04-08"
$ws.Cells.Item(115, 3).Value = 0.239
$ws.Cells.Item(116, 1).Value = "2024-04-08"
$ws.Cells.Item(116, 3).Value = 0.376
$ws.Cells.Item(117, 1).Value = "2024-04-08"
$ws.Cells.Item(117, 3).Value = 0.488
$ws.Cells.Item(118, 1).Value = "2024-04-08"
$ws.Cells.Item(118, 3).Value = 0.548
$ws.Cells.Item(119, 1).Value = "2024-04-08"
$ws.Cells.Item(119, 3).Value = 0.564
$ws.Cells.Item(120, 1).Value = "2024-04-08"
$ws.Cells.Item(120, 3).Value = 0.539
$ws.Cells.Item(121, 1).Value = "2024-04-08"
$ws.Cells.Item(121, 3).Value = 0.512
$ws.Cells.Item(122, 1).Value = "2024-04-08"
$ws.Cells.Item(122, 3).Value = 0.416
$ws.Cells.Item(123, 1).Value = "2024-04-08"
$ws.Cells.Item(123, 3).Value = 0.303
$ws.Cells.Item(124, 1).Value = "2024-04-08"
$ws.Cells.Item(124, 3).Value = 0.159
$ws.Cells.Item(125, 1).Value = "2024-04-08"
$ws.Cells.Item(125, 3).Value = 0.03
$ws.Cells.Item(126, 1).Value = "2024-04-08"
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(127, 1).Value = "2024-04-08"
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(128, 1).Value = "2024-04-08"
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(129, 1).Value = "2024-04-08"
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(130, 1).Value = "2024-04-09"
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(131, 1).Value = "2024-04-09"
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(132, 1).Value = "2024-04-09"
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(133, 1).Value = "2024-04-09"
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(134, 1).Value = "2024-04-09"
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(135, 1).Value = "2024-04-09"
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(136, 1).Value = "2024-04-09"
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(137, 1).Value = "2024-04-09"
$ws.Cells.Item(137, 3).Value = 0.024
$ws.Cells.Item(138, 1).Value = "2024-04-09"
$ws.Cells.Item(138, 3).Value = 0.116
$ws.Cells.Item(139, 1).Value = "2024-04-09"
$ws.Cells.Item(139, 3).Value = 0.242
$ws.Cells.Item(140, 1).Value = "2024-04-09"
$ws.Cells.Item(140, 3).Value = 0.383
$ws.Cells.Item(141, 1).Value = "2024-04-09"
$ws.Cells.Item(141, 3).Value = 0.506
$ws.Cells.Item(142, 1).Value = "2024-04-09"
$ws.Cells.Item(142, 3).Value = 0.552
$ws.Cells.Item(143, 1).Value = "2024-04-09"
$ws.Cells.Item(143, 3).Value = 0.564
$ws.Cells.Item(144, 1).Value = "2024-04-09"
$ws.Cells.Item(144, 3).Value = 0.539
$ws.Cells.Item(145, 1).Value = "2024-04-09"
$ws.Cells.Item(145, 3).Value = 0.486
$ws.Cells.Item(146, 1).Value = "2024-04-09"
$ws.Cells.Item(146, 3).Value = 0.385
$ws.Cells.Item(147, 1).Value = "2024-04-09"
$ws.Cells.Item(147, 3).Value = 0.274
$ws.Cells.Item(148, 1).Value = "2024-04-09"
$ws.Cells.Item(148, 3).Value = 0.162
$ws.Cells.Item(149, 1).Value = "2024-04-09"
$ws.Cells.Item(149, 3).Value = 0.028
$ws.Cells.Item(150, 1).Value = "2024-04-09"
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(151, 1).Value = "2024-04-09"
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(152, 1).Value = "2024-04-09"
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(153, 1).Value = "2024-04-09"
$ws.Cells.Item(153, 3).Value = 0
$ws.Cells.Item(154, 1).Value = "2024-04-10"
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(155, 1).Value = "2024-04-10"
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(156, 1).Value = "2024-04-10"
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(157, 1).Value = "2024-04-10"
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(158, 1).Value = "2024-04-10"
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(159, 1).Value = "2024-04-10"
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(160, 1).Value = "2024-04-10"
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(161, 1).Value = "2024-04-10"
$ws.Cells.Item(161, 3).Value = 0.028
$ws.Cells.Item(162, 1).Value = "2024-04-10"
$ws.Cells.Item(162, 3).Value = 0.129
$ws.Cells.Item(163, 1).Value = "2024-04-10"
$ws.Cells.Item(163, 3).Value = 0.268
$ws.Cells.Item(164, 1).Value = "2024-04-10"
$ws.Cells.Item(164, 3).Value = 0.402
$ws.Cells.Item(165, 1).Value = "2024-04-10"
$ws.Cells.Item(165, 3).Value = 0.506
$ws.Cells.Item(166, 1).Value = "2024-04-10"
$ws.Cells.Item(166, 3).Value = 0.563
$ws.Cells.Item(167, 1).Value = "2024-04-10"
$ws.Cells.Item(167, 3).Value = 0.573
$ws.Cells.Item(168, 1).Value = "2024-04-10"
$ws.Cells.Item(168, 3).Value = 0.565
$ws.Cells.Item(169, 1).Value = "2024-04-10"
$ws.Cells.Item(169, 3).Value = 0.499
$ws.Cells.Item(170, 1).Value = "2024-04-10"
$ws.Cells.Item(170, 3).Value = 0.379

# Restore the original date display format (values remain text, as in
# the source file, since they were authored as inline strings).
$ws.Range("A2:A170").NumberFormat = "dd.mm.yyyy"
